# Apply the "for loops in for brs, aaq, flourish" edit to sheet1.
# Net effect vs before.xlsx:
#  - Row 8 (anx_sev) is removed from column A, shifting the variable-name
#    labels in A9:A63 up by one row (B/C values stay put per row).
#  - A couple of labels are renamed in place: aaq_yn -> aaq_dum,
#    anx_mod -> anx_score, BRS_tot -> brs_dum, flourish -> flourish_dum.
#  - The now-empty last row (63) is removed, shrinking the used range to
#    A1:C62.
# We simply rewrite column A for rows 2-62 to the final expected text and
# clear out row 63 completely (which also updates the sheet dimension).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value2  = "aaq_dum"
$ws.Range("A7").Value2  = "anx_score"
$ws.Range("A8").Value2  = "assault_emo"
$ws.Range("A9").Value2  = "assault_phys"
$ws.Range("A10").Value2 = "assault_sex"
$ws.Range("A11").Value2 = "audit_tot"
$ws.Range("A12").Value2 = "belong1"
$ws.Range("A13").Value2 = "belong2"
$ws.Range("A14").Value2 = "belong8"
$ws.Range("A15").Value2 = "belong9"
$ws.Range("A16").Value2 = "binge_fr"
$ws.Range("A17").Value2 = "body_sr"
$ws.Range("A18").Value2 = "brs_dum"
$ws.Range("A19").Value2 = "dep_impa"
$ws.Range("A20").Value2 = "dep_secret"
$ws.Range("A21").Value2 = "deprawsc"
$ws.Range("A22").Value2 = "discrim"
$ws.Range("A23").Value2 = "divers"
$ws.Range("A24").Value2 = "drug_mar"
$ws.Range("A25").Value2 = "drugs_yn"
$ws.Range("A26").Value2 = "dx_adhd"
$ws.Range("A27").Value2 = "dx_bi"
$ws.Range("A28").Value2 = "dx_dep"
$ws.Range("A29").Value2 = "dx_pers"
$ws.Range("A30").Value2 = "dx_tr"
$ws.Range("A31").Value2 = "ed_any"
$ws.Range("A32").Value2 = "env_mh"
$ws.Range("A33").Value2 = "fincur"
$ws.Range("A34").Value2 = "finpast"
$ws.Range("A35").Value2 = "flourish_dum"
$ws.Range("A36").Value2 = "gad7_impa"
$ws.Range("A37").Value2 = "gender_noncis"
$ws.Range("A38").Value2 = "gpa_sr"
$ws.Range("A39").Value2 = "inf"
$ws.Range("A40").Value2 = "ins_cover"
$ws.Range("A41").Value2 = "international"
$ws.Range("A42").Value2 = "meds_anx"
$ws.Range("A43").Value2 = "meds_count"
$ws.Range("A44").Value2 = "meds_dep"
$ws.Range("A45").Value2 = "meds_mood"
$ws.Range("A46").Value2 = "meds_sle"
$ws.Range("A47").Value2 = "meds_sti"
$ws.Range("A48").Value2 = "military"
$ws.Range("A49").Value2 = "percneed_cur"
$ws.Range("A50").Value2 = "persist"
$ws.Range("A51").Value2 = "psyhx"
$ws.Range("A52").Value2 = "race"
$ws.Range("A53").Value2 = "religios"
$ws.Range("A54").Value2 = "residenc"
$ws.Range("A55").Value2 = "satisfied_overall"
$ws.Range("A56").Value2 = "school2_type"
$ws.Range("A57").Value2 = "sexual"
$ws.Range("A58").Value2 = "sib_freq"
$ws.Range("A59").Value2 = "stig_pcv_2"
$ws.Range("A60").Value2 = "stig_pcv_3"
$ws.Range("A61").Value2 = "talk"
$ws.Range("A62").Value2 = "ther_vis"

# Remove the now-redundant last row entirely (data, not just column A),
# which also shrinks the sheet dimension down to A1:C62.
$ws.Range("A63:C63").ClearContents()
